$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage for cells whose new
# value would otherwise be auto-parsed by Excel as a number (General format).
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "29.104.80"
$ws.Range("D3").Value = "1.843.74"
$ws.Range("E3").Value = "  -2.10%  "
Set-TextValue "D4" "1.0000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -5.94%  "
Set-TextValue "D6" "237.56"
$ws.Range("E6").Value = "  -2.26%  "
Set-TextValue "D7" "0.9998"
Set-TextValue "D8" "0.3031"
$ws.Range("E8").Value = "  -4.11%  "
Set-TextValue "D9" "0.07435"
$ws.Range("E9").Value = "  +2.98%  "
Set-TextValue "D10" "23.26"
$ws.Range("E10").Value = "  -6.69%  "
Set-TextValue "D11" "0.08100"
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D12" "0.7226"
$ws.Range("E12").Value = "  -4.44%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "5.226"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.816.38"
$ws.Range("E14").Value = "  -5.43%  "
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("D16").Value = "29.093.35"
$ws.Range("E16").Value = "  -2.83%  "
Set-TextValue "D17" "5.788"
$ws.Range("E17").Value = "  -6.25%  "
Set-TextValue "D18" "240.41"
$ws.Range("E18").Value = "  -3.79%  "
Set-TextValue "D19" "0.000007654"
$ws.Range("E19").Value = "  -2.54%  "
Set-TextValue "D20" "13.00"
$ws.Range("E20").Value = "  -4.54%  "
Set-TextValue "D21" "0.9986"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "2.086.63"
$ws.Range("E22").Value = "  -2.19%  "
Set-TextValue "D23" "1.000"
$ws.Range("E23").Value = "  +0.03%  "
Set-TextValue "D24" "7.550"
$ws.Range("E24").Value = "  -6.27%  "
Set-TextValue "D25" "161.83"
$ws.Range("E25").Value = "  -2.23%  "
Set-TextValue "D26" "0.1459"
$ws.Range("E26").Value = "  -6.56%  "
Set-TextValue "D27" "8.937"
$ws.Range("E27").Value = "  -3.94%  "
Set-TextValue "D28" "17.96"
$ws.Range("E28").Value = "  -4.12%  "
$ws.Range("E29").Value = "  -5.67%  "
$ws.Range("E30").Value = "  -8.15%  "
Set-TextValue "D31" "4.473"
$ws.Range("E31").Value = "  -3.11%  "
Set-TextValue "D32" "1.490"
$ws.Range("E32").Value = "  -3.09%  "
Set-TextValue "D33" "4.015"
Set-TextValue "D34" "0.05171"
$ws.Range("E34").Value = "  -3.77%  "
Set-TextValue "D35" "1.182"
$ws.Range("E35").Value = "  -5.84%  "
Set-TextValue "D36" "0.7090"
$ws.Range("E36").Value = "  -6.67%  "
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  -2.18%  "
Set-TextValue "D39" "0.01867"
$ws.Range("E39").Value = "  -5.15%  "
Set-TextValue "D40" "2.670"
$ws.Range("E40").Value = "  -3.37%  "
Set-TextValue "D41" "0.8998"
$ws.Range("E41").Value = "  +3.22%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "5.919"
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D43" "0.4272"
$ws.Range("E43").Value = "  -6.41%  "
$ws.Range("D44").Value = "1.054.22"
$ws.Range("E44").Value = "  -4.62%  "
Set-TextValue "D45" "69.71"
$ws.Range("E45").Value = "  -4.01%  "
Set-TextValue "D46" "0.9996"
$ws.Range("E46").Value = "  -0.09%  "
Set-TextValue "D47" "101.45"
$ws.Range("E47").Value = "  -3.14%  "
Set-TextValue "D48" "1.748"
$ws.Range("E48").Value = "  -6.50%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.986.87"
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "9.186"
$ws.Range("E50").Value = "  -3.88%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D51" "7.038"
$ws.Range("E51").Value = "  -7.62%  "
